# Adds a new "2022-Q3" worksheet (with its single fund-holding row) right
# after "总计" and before the former first quarter sheet, then records that
# quarter in the "总计" summary sheet as its new row 2 (pushing the rest of
# the history down by one row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1 — insert the new "2022-Q3" sheet, positioned right after "总计".
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)          # "总计"
$oldFirstQuarter = $wb.Worksheets.Item(2)      # currently "2022-Q2"
$newSheet = $wb.Worksheets.Add($oldFirstQuarter)
$newSheet.Name = "2022-Q3"

# Re-use the column layout/formatting (bold + bordered header row, bordered
# index column) from the sheet that used to be first — now shifted to slot 3.
$q2Sheet = $wb.Worksheets.Item(3)
$q2Sheet.Range("B1:H2").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)      # xlPasteFormats
$q2Sheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)      # xlPasteFormats

# Header row.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Single holding row. Fund code / size / position columns are text in every
# quarter sheet (leading zeros such as "010714" must survive), so they are
# entered with a leading apostrophe to stop Excel auto-coercing them to
# numbers, then the quote-prefix styling left behind is cleared again.
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'010714"
$newSheet.Range("C2").Value = "东方红远见价值混合"
$newSheet.Range("D2").Value = "'15.24"
$newSheet.Range("E2").Value = "'94.15"
$newSheet.Range("F2").Value = "'3.04"
$newSheet.Range("G2").Value = "'0.4633"
$newSheet.Range("H2").Value = 10
$newSheet.Range("B2:G2").Style = "Normal"

# ---------------------------------------------------------------------------
# Step 2 — record the new quarter as row 2 of the "总计" summary sheet,
# shifting the existing history rows down by one.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.46

# The row insert leaves the new row's numeric-looking cells without the
# workbook's usual formatting; copy it over explicitly from the row below
# (which already carries the correct styles: bordered/bold index cell,
# plain data cells) so the result matches the rest of the table.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)    # xlPasteFormats
$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122) # xlPasteFormats

Write-Output "2022-Q3 sheet added and 总计 summary updated"
